$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 0. Drop every existing hyperlink up front. Several of them need to move
#    to a different cell and the host engine does not re-anchor hyperlink
#    refs when rows are edited/deleted, so the safest approach is to
#    remove them all now and recreate the final, correct set at the end.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 1. Row 13 ("Hack-E-Thon" / "Roll-Off-Form") becomes a section header:
#    give B13 the same bold "header" formatting already used by
#    B1 / B6 / B10 (value itself is unchanged).
# ---------------------------------------------------------------------
$ws.Range("B1").Copy()
$ws.Range("B13").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. New row 14: "Frontend and Backend" / "Link"
# ---------------------------------------------------------------------
$ws.Range("A14").Value = "Frontend and Backend"
$ws.Range("A2").Copy()
$ws.Range("A14").PasteSpecial(-4122)

$ws.Range("B14").Value = "Link"
$ws.Range("B2").Copy()
$ws.Range("B14").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Row 15 loses its "Assignments(Training)" label (moves to row 19).
# ---------------------------------------------------------------------
$ws.Range("A15").Value = ""

# ---------------------------------------------------------------------
# 4. Row 16: "One Pager Resume in CapGemini Format" / "Link"
#    (this used to live on row 17)
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "One Pager Resume in CapGemini Format"
$ws.Range("A2").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("B16").Value = "Link"
$ws.Range("B2").Copy()
$ws.Range("B16").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 5. Row 17: "Video link of your Case Study" / "Link"
#    (this used to live on row 18). A17/B17 already carry the right
#    cell styles, only the text needs touching.
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "Video link of your Case Study"
$ws.Range("B17").Value = "Link"

# ---------------------------------------------------------------------
# 6. Row 18 becomes fully blank (A18 text cleared, B18 cell removed).
# ---------------------------------------------------------------------
$ws.Range("A18").Value = ""
$ws.Range("B18").Clear()

# ---------------------------------------------------------------------
# 7. Row 19: "Assignments(Training)" / "Link" (new link row).
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "Assignments(Training)"
$ws.Range("A2").Copy()
$ws.Range("A19").PasteSpecial(-4122)

$ws.Range("B19").Value = "Link"
$ws.Range("B2").Copy()
$ws.Range("B19").PasteSpecial(-4122)

# Row 20 ("Post training ->CAF Completion Screen shot" / "Link") is untouched.

# ---------------------------------------------------------------------
# 8. Row 21: "Certification Updates -> Screen shot and r2d2 upload" / "Link"
#    (this used to live on row 24).
# ---------------------------------------------------------------------
$ws.Range("A21").Value = "Certification Updates -> Screen shot and r2d2 upload"
$ws.Range("A2").Copy()
$ws.Range("A21").PasteSpecial(-4122)

$ws.Range("B21").Value = "Link"
$ws.Range("B2").Copy()
$ws.Range("B21").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 9. Drop the now-unused trailing rows (23, 24, 25), bottom row first so
#    the indices of the rows still to be removed don't shift underneath us.
# ---------------------------------------------------------------------
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(23).Delete()

# ---------------------------------------------------------------------
# 10. Recreate every hyperlink against its final cell.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/shubhamrajj/RailwayReservationSystem.git", "", "", "Link")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/shubhamrajj/Railway-Reservation-System-Server.git", "", "", "Link")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/shubhamrajj/RailwayReservationSystemDocument.git", "", "", "Link")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/shubhamrajj/BankingSystem_client.git", "", "", "Link")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://github.com/shubhamrajj/BankingProjectAPI.git", "", "", "Link")
$ws.Hyperlinks.Add($ws.Range("B11"), "https://github.com/shubhamrajj/Updated_EmployeePerformance-master.git", "", "", "Link")
$ws.Hyperlinks.Add($ws.Range("B14"), "https://github.com/shubhamrajj/Hack-E-Thon.git", "", "", "Link")
$ws.Hyperlinks.Add($ws.Range("B16"), "https://github.com/shubhamrajj/Resume.git", "", "", "Link")
$ws.Hyperlinks.Add($ws.Range("B17"), "https://www.loom.com/share/094359b7be1d4238b641b9594d3c760c", "", "", "Link")
$ws.Hyperlinks.Add($ws.Range("B19"), "https://github.com/shubhamrajj/Assignments.git", "", "", "Link")
$ws.Hyperlinks.Add($ws.Range("B20"), "https://github.com/shubhamrajj/CAF.git", "", "", "Link")
$ws.Hyperlinks.Add($ws.Range("B21"), "https://github.com/shubhamrajj/Certificates.git", "", "", "Link")

# Re-apply the normal hyperlink cell look (Hyperlinks.Add stamps a fresh
# style the first time it is used) so every link cell still shares the
# same cell style as the rest of the sheet.
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("B2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 11. Selection moves to F18.
# ---------------------------------------------------------------------
$ws.Range("F18").Select()

Write-Host "Edits applied"
